$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, E, F hold numeric-looking identifiers that must remain text,
# matching the inlineStr type used in the source workbook.
$ws.Range("D2:F4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '281474991265526-1748559988065'
$ws.Range("B2").Value = 'Camera Obstruction'
$ws.Range("C2").Value = '2025-05-29T17:06:28.065'
$ws.Range("D2").Value = '281474991265526'
$ws.Range("E2").Value = '150'
$ws.Range("F2").Value = '51834047'
$ws.Range("G2").Value = 'EDUARDO HERNANDEZ'
$ws.Range("H2").Value = 20.36509002
$ws.Range("I2").Value = -103.20564524
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 'No video URL'
$ws.Range("L2").Value = 'No video URL'

# Row 3
$ws.Range("A3").Value = '281474991395097-1748544268163'
$ws.Range("B3").Value = 'Harsh Turn'
$ws.Range("C3").Value = '2025-05-29T12:44:28.163'
$ws.Range("D3").Value = '281474991395097'
$ws.Range("E3").Value = '125'
$ws.Range("F3").Value = '51834055'
$ws.Range("G3").Value = 'DAVID SERRANO'
$ws.Range("H3").Value = 20.651197249
$ws.Range("I3").Value = -103.353378809
$ws.Range("J3").Value = 0.6844108700752258
$ws.Range("K3").Value = 'No video URL'
$ws.Range("L3").Value = 'No video URL'

# Row 4
$ws.Range("A4").Value = '281474991452844-1748528546412'
$ws.Range("B4").Value = 'No Seat Belt'
$ws.Range("C4").Value = '2025-05-29T08:22:26.412'
$ws.Range("D4").Value = '281474991452844'
$ws.Range("E4").Value = '134'
$ws.Range("F4").Value = '51834052'
$ws.Range("G4").Value = 'DIEGO CÁRDENAS'
$ws.Range("H4").Value = 20.64218943
$ws.Range("I4").Value = -103.34121831
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991452844/1748528543912/W2Y6Fn2i0v-camera-video-segment-driver-1748528546412.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSL6X4UFNU%2F20250530%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250530T160221Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjENv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDuv30BgBOWrrv00vRVivFCbUJ5tSotHrtb4CliaX2RwAIgKgZsDBaac8dBbr0iv8FVCrzUdcQ8zioTu7OfHpziNC4q5gMIpP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEOfACgm%2F6yIvn4ocyq6A0%2FB1PNVvwcEiFYNURKzD2ApRUXZbr5yzS6Pb3TrwfBEx7EAg5W6xGq42qmXXD9%2FJ98BGJEnxf8vnRXJLzQeBnsoGOYoaG%2BQH3nuiP%2F3KFg0iQBC89E%2BVWAPMVIg30Z65Cig3l6tmqyzufFlRO%2Fd3PP%2FAylb1IelSEG94oKstwQFjzfB%2F78zv4fQrgMOBjLQo9T5lV3HURS4zCh%2Fm90ygtk%2FOKTo8Wu%2F%2FFCkvORVZ4lDin3vTgOP4mjDUAeNC%2F67lK76Dur%2BX9RLYRTlfAIBv7qvX5BKvTD4p%2FTRCpmDN84p2wdRFSHKlcrMGza%2Fr%2Br9JOvdUylpwqCdgLiIs%2BuTUWxarMEBrhcm2SBSQzlRykCGjvHPgj0ELfE52svJ9gjlNb32bhWC5ijpDM43zsIBr9Lb%2Bum7eHraqqx74WzGwBpHv%2BcxAP094aPCBfDmbPN88G3C81kw00WdBeUt27lc9XuuoVGmi5U2ByxCf6%2BPeph30IBrrnDRLY4Bw28E239soNH6TPZ0cxrUzPIDMaIziAMTK8j01FTfChcRCvaby1EmRsY67X9f4YzWwhRP0C5k9ZefeRd%2BDBzK6dkwp6fmwQY6pQHryzeVn%2Ft2fikFSoDr04UB5vPN2LFrqIhwv%2FPqNpQkOKXq58T74Le%2BGglgufostYqzWHRt%2BsSl9wjIpNRg8qM5JejYf5A9QxnGD2%2FVl1g3JpmIOk%2FJphE9wg7GCoTIrc%2FkRG%2BHVrhPln2osjnxhshj2Kf29qXID7vZKg96M2keJ%2FgD4VyCZUZtM5NwApRaYgZad7V702im3RZBo60SbchhKVbACEo%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2000%3A02%3A21%20GMT&X-Amz-Signature=a612deef59e1beae5c4a5402124e39df8ef0a04ddf3da370655b25bbb56bf2b7'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991452844/1748528543912/rQWZBAyQY2-camera-video-segment-1748528546412.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSL6X4UFNU%2F20250530%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250530T160221Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjENv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDuv30BgBOWrrv00vRVivFCbUJ5tSotHrtb4CliaX2RwAIgKgZsDBaac8dBbr0iv8FVCrzUdcQ8zioTu7OfHpziNC4q5gMIpP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEOfACgm%2F6yIvn4ocyq6A0%2FB1PNVvwcEiFYNURKzD2ApRUXZbr5yzS6Pb3TrwfBEx7EAg5W6xGq42qmXXD9%2FJ98BGJEnxf8vnRXJLzQeBnsoGOYoaG%2BQH3nuiP%2F3KFg0iQBC89E%2BVWAPMVIg30Z65Cig3l6tmqyzufFlRO%2Fd3PP%2FAylb1IelSEG94oKstwQFjzfB%2F78zv4fQrgMOBjLQo9T5lV3HURS4zCh%2Fm90ygtk%2FOKTo8Wu%2F%2FFCkvORVZ4lDin3vTgOP4mjDUAeNC%2F67lK76Dur%2BX9RLYRTlfAIBv7qvX5BKvTD4p%2FTRCpmDN84p2wdRFSHKlcrMGza%2Fr%2Br9JOvdUylpwqCdgLiIs%2BuTUWxarMEBrhcm2SBSQzlRykCGjvHPgj0ELfE52svJ9gjlNb32bhWC5ijpDM43zsIBr9Lb%2Bum7eHraqqx74WzGwBpHv%2BcxAP094aPCBfDmbPN88G3C81kw00WdBeUt27lc9XuuoVGmi5U2ByxCf6%2BPeph30IBrrnDRLY4Bw28E239soNH6TPZ0cxrUzPIDMaIziAMTK8j01FTfChcRCvaby1EmRsY67X9f4YzWwhRP0C5k9ZefeRd%2BDBzK6dkwp6fmwQY6pQHryzeVn%2Ft2fikFSoDr04UB5vPN2LFrqIhwv%2FPqNpQkOKXq58T74Le%2BGglgufostYqzWHRt%2BsSl9wjIpNRg8qM5JejYf5A9QxnGD2%2FVl1g3JpmIOk%2FJphE9wg7GCoTIrc%2FkRG%2BHVrhPln2osjnxhshj2Kf29qXID7vZKg96M2keJ%2FgD4VyCZUZtM5NwApRaYgZad7V702im3RZBo60SbchhKVbACEo%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2031%20May%202025%2000%3A02%3A21%20GMT&X-Amz-Signature=ec59af9e0bed2a3b3ffecfb146afd924d8eecfc0e2d5019b6a9e6024fc516e7a'

